$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "https://leetcode.com/problems/maximum-subarray/description/"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "Guided"
$ws.Range("D5").Value = "Guided"
$ws.Range("E5").Value = "Self"
